$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 25.02.2022 12:00"

# Row 10 (EuroOil Opuštěná): convert D10/E10 from text to real numeric values
# D10: delta price as a number (was text "+0.2")
$ws.Range("D10").Value = 0.2

# E10: date/time as an Excel serial number, matching the format used by the
# other rows in column E (numFmt "YYYY-MM-DD HH:MM:SS"). Set the number
# format first so Excel doesn't create an extra implicit date style.
$ws.Range("E10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = (Get-Date -Year 2022 -Month 2 -Day 25 -Hour 11 -Minute 47 -Second 32)
